$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab from "Sheet1" to "Array"
$ws.Name = "Array"

# 2. Insert 4 new rows right after the existing MPG "<40" row (old row 16) to hold the
#    new "<45", "<50", "<55", "<60" MPG buckets. This pushes everything from the old
#    row 17 ("sedan") onward down by 4 rows.
$ws.Range("A17:A20").EntireRow.Insert()

# Fill in the 4 new MPG rows and give them the same fill/style as the rest of the
# MPG block (copy format from B16, which carries style index 1).
$ws.Range("B17").Value = "<45"
$ws.Range("B18").Value = "<50"
$ws.Range("B19").Value = "<55"
$ws.Range("B20").Value = "<60"
$ws.Range("B16").Copy()
$ws.Range("B17:B20").PasteSpecial(-4122)

# Move the "MPG" category label so it stays centred on the now-9-row MPG block:
# it used to live on C14 (old block 12-16), now it belongs on C16 (new block 12-20).
$ws.Range("C14").Clear()
$ws.Range("C16").Value = "MPG"
$ws.Range("B16").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# 3. Remove the "electric" row from the Engine Fuel Type block. After the insert above,
#    the old row 28 ("electric") is now row 32.
$ws.Range("A32").EntireRow.Delete()

# Move the "Engine Fuel Type" label so it stays centred on the now-3-row block
# (gas / hybrid / diesel): it belongs on row 32 (the "hybrid" row) in columns C and D,
# matching the block's style (index 2).
$ws.Range("C32").Value = "Engine Fuel Type"
$ws.Range("B32").Copy()
$ws.Range("C32:D32").PasteSpecial(-4122)

# 4. Renumber the sequential index column (A3:A49) back to a clean 0..46 run, since the
#    row insert/delete above left gaps/duplicates in the numbering.
for ($i = 0; $i -le 46; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $i
}

# 5. Sheet view: drop the old scrolled/selected state and set the new selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E15").Select()
